$p = $ppt.ActivePresentation
$hm = $p.HandoutMaster
try {
  $tcs = $hm.ThemeColorScheme
  Write-Host "Handout TCS: $tcs Count=$($tcs.Count)"
} catch { Write-Host "err: $_" }
